$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("K2:K54").Value = "pb"
